$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add three new "Pre-irrad" loading-diagram sheets (9, 10, 11) by copying the
# previous sheet in the chain each time and repointing its formulas at the
# sheet it was copied from (mirrors how every earlier Pre-irrad_N sheet in
# this workbook was produced).
# ---------------------------------------------------------------------------

$chain = @(
    @{ Prev = "Pre-irrad_8_141107";  New = "Pre-irrad_9_141110";  Cell = "A1" },
    @{ Prev = "Pre-irrad_9_141110";  New = "Pre-irrad_10_141113"; Cell = "A2" },
    @{ Prev = "Pre-irrad_10_141113"; New = "Pre-irrad_11_141114"; Cell = "A1" }
)

# "Pre-irrad_8_141107"'s own formulas still reference the sheet IT was
# copied from ("Pre-irrad_7_141106") -- that's the text actually present in
# the copied cells, so that's what the first Replace must search for.
$referencedSheet = "Pre-irrad_7_141106"

foreach ($step in $chain) {
    $src = $wb.Worksheets.Item($step.Prev)
    $src.Copy($null, $src)
    $new = $wb.Worksheets.Item($wb.Worksheets.Count)
    $new.Name = $step.New

    # The copied formulas still point at $referencedSheet; repoint them at
    # $src itself (the sheet this new one was actually copied from).
    $new.UsedRange.Replace($referencedSheet, $step.Prev) | Out-Null

    # Reset the view: plain A1 (or A2) selection, tab not selected yet.
    $new.Range($step.Cell).Select() | Out-Null

    # Now $new's formulas reference $step.Prev, so that's the text the
    # *next* sheet in the chain will need repointed.
    $referencedSheet = $step.Prev
}

# The previously-last sheet ("Pre-irrad_8_141107") is no longer the active
# tab or selected at A48 -- it settles back to A1, unselected.
$ws8 = $wb.Worksheets.Item("Pre-irrad_8_141107")
$ws8.Range("A1").Select() | Out-Null

# "Pre-irrad_6_141105" view resets from C34 back to A1 as well.
$ws6 = $wb.Worksheets.Item("Pre-irrad_6_141105")
$ws6.Range("A1").Select() | Out-Null

# The newest sheet ("Pre-irrad_11_141114") becomes the active tab.
$wsLast = $wb.Worksheets.Item("Pre-irrad_11_141114")
$wsLast.Activate()
$wsLast.Range("A1").Select() | Out-Null
